$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: volume number and week-covering dates ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Crime Complaints table (rows 14-30): new weekly data ---
$ws.Range("D14").Value = 3
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").Value = 4
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = -20
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = -22.222222222222
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = -30
$ws.Range("N15").Value = -70.833333333333
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -35
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 137
$ws.Range("K16").Value = -10.948905109489
$ws.Range("L16").Value = 74.285714285714
$ws.Range("M16").Value = -37.755102040816
$ws.Range("N16").Value = -84.258064516129
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 175
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 2.941176470588
$ws.Range("I17").Value = 271
$ws.Range("J17").Value = 280
$ws.Range("K17").Value = -3.214285714285
$ws.Range("L17").Value = 7.539682539682
$ws.Range("M17").Value = 69.375
$ws.Range("N17").Value = -52.204585537918
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -41.176470588235
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = -31.578947368421
$ws.Range("L18").Value = -37.6
$ws.Range("M18").Value = -35
$ws.Range("N18").Value = -89.090909090909
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 37.5
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -35.555555555555
$ws.Range("I19").Value = 346
$ws.Range("J19").Value = 378
$ws.Range("K19").Value = -8.465608465608
$ws.Range("L19").Value = 23.571428571428
$ws.Range("M19").Value = -21.719457013574
$ws.Range("N19").Value = -29.95951417004
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -7.142857142857
$ws.Range("I20").Value = 78
$ws.Range("J20").Value = 91
$ws.Range("K20").Value = -14.285714285714
$ws.Range("L20").Value = 52.941176470588
$ws.Range("M20").Value = -7.142857142857
$ws.Range("N20").Value = -90.681003584229
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -3.448275862068
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = -24.444444444444
$ws.Range("I21").Value = 906
$ws.Range("J21").Value = 1014
$ws.Range("K21").Value = -10.650887573964
$ws.Range("L21").Value = 14.393939393939
$ws.Range("M21").Value = -11.001964636542
$ws.Range("N21").Value = -73.586005830903
$ws.Range("G22").Value = 2
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = -53.846153846153
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 500
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 122
$ws.Range("J23").Value = 116
$ws.Range("K23").Value = 5.172413793103
$ws.Range("L23").Value = -8.270676691729
$ws.Range("M23").Value = 40.229885057471
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -39.285714285714
$ws.Range("F24").Value = 71
$ws.Range("H24").Value = -28.282828282828
$ws.Range("I24").Value = 740
$ws.Range("J24").Value = 886
$ws.Range("K24").Value = -16.47855530474
$ws.Range("L24").Value = 24.369747899159
$ws.Range("M24").Value = -19.302071973827
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -46.666666666666
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 449
$ws.Range("J25").Value = 418
$ws.Range("K25").Value = 7.416267942583
$ws.Range("L25").Value = 40.752351097178
$ws.Range("M25").Value = 16.020671834625
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 13
$ws.Range("K26").Value = -23.529411764705
$ws.Range("L26").Value = -23.529411764705
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 40
$ws.Range("K27").Value = 8.108108108108
$ws.Range("L27").Value = 53.846153846153
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -85.714285714285
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 26
$ws.Range("K28").Value = -73.076923076923
$ws.Range("L28").Value = -30
$ws.Range("M28").Value = -53.333333333333
$ws.Range("N28").Value = -85.416666666666
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 17
$ws.Range("K29").Value = -76.470588235294
$ws.Range("L29").Value = -60
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -89.473684210526